# Dev Guide diagram refresh:
#  - bump every cached "datetimeFigureOut" field (notes master, slide
#    master, all custom layouts) from 1/7/2017 to 11/11/18
#  - rename the BrowserPanel rectangle on slide 1 to StaffPanel

$p = $ppt.ActivePresentation

$oldDate = "1/7/2017"
$newDate = "11/11/18"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder.
$sm = $p.SlideMaster
Update-DateShape $sm.Shapes

# Every custom layout hanging off the slide master.
$layouts = $sm.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $cl = $layouts.Item($li)
    Update-DateShape $cl.Shapes
}

# Notes master date placeholder - only reachable through the
# HeadersFooters.DateAndTime sub-object on this host.
$nm = $p.NotesMaster
$nmHf = $nm.HeadersFooters
$nmHf.DateAndTime.Text = $newDate

# Slide 1: rename the BrowserPanel rectangle to StaffPanel.
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "BrowserPanel") {
            $sh.TextFrame.TextRange.Text = "StaffPanel"
        }
    }
}
